$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit inserts a new company record ("江苏长江汇科技有限公司") as the new
# row 2 of the data table. All subsequent rows (previously rows 2-15) shift
# down by one row (now rows 3-16). The numbering in column A stays pinned to
# each physical row position (0,1,2,...), so a brand new value (14) is appended
# for the newly created last row (16), which now holds what used to be row 15.
# ---------------------------------------------------------------------------

# ---- Row 2 ----
$ws.Cells.Item(2, 2).Value = '江苏长江汇科技有限公司'
$ws.Cells.Item(2, 3).Value = '鼓楼区'
$ws.Cells.Item(2, 4).Value = '技术部'
$ws.Cells.Item(2, 5).Value = 'Java'
$ws.Cells.Item(2, 6).Value = '9:00-17:30'
$ws.Cells.Item(2, 7).Value = '1.5h'
$ws.Cells.Item(2, 8).Value = '按需加班'
$ws.Cells.Item(2, 9).Value = '基数 看个人，比例 5%'
$ws.Cells.Item(2, 10).Value = '基数5000 按照绩效或多或少'
$ws.Cells.Item(2, 11).Value = '8折'
$ws.Cells.Item(2, 12).Value = '无隔板工位，配台式电脑，自带电脑每个月有200补贴，持续24个月。'
$ws.Cells.Item(2, 13).Value = ''
$ws.Cells.Item(2, 14).Value = '钉钉严格打卡'
$ws.Cells.Item(2, 15).Value = ''
$ws.Cells.Item(2, 16).Value = ''
$ws.Cells.Item(2, 17).Value = '2022-02-06 13:30:06'

# ---- Row 3 ----
$ws.Cells.Item(3, 2).Value = '南京伯索网络科技有限公司（PLASO）'
$ws.Cells.Item(3, 3).Value = '秦淮区'
$ws.Cells.Item(3, 4).Value = ''
$ws.Cells.Item(3, 5).Value = ''
$ws.Cells.Item(3, 6).Value = '9:00-18:00'
$ws.Cells.Item(3, 7).Value = '1h'
$ws.Cells.Item(3, 8).Value = '124 加班，35 正常；大小周'
$ws.Cells.Item(3, 9).Value = '基数南京底薪，比例 8%'
$ws.Cells.Item(3, 10).Value = '一般无'
$ws.Cells.Item(3, 11).Value = '3个月8折'
$ws.Cells.Item(3, 12).Value = '网吧工位'
$ws.Cells.Item(3, 13).Value = '入职一年后才有，每年加一天'
$ws.Cells.Item(3, 14).Value = '企业微信打卡，每月三次迟到机会'
$ws.Cells.Item(3, 15).Value = ''
$ws.Cells.Item(3, 16).Value = ''
$ws.Cells.Item(3, 17).Value = '2022-02-06 13:26:16'

# ---- Row 4 ----
$ws.Cells.Item(4, 2).Value = '南京叶子科技有限公司'
$ws.Cells.Item(4, 3).Value = 'xx区'
$ws.Cells.Item(4, 4).Value = 'xxx事业部'
$ws.Cells.Item(4, 5).Value = 'Java'
$ws.Cells.Item(4, 6).Value = '9:00-18:30'
$ws.Cells.Item(4, 7).Value = '1.5h'
$ws.Cells.Item(4, 8).Value = '135 加班，24 正常；大小周等等'
$ws.Cells.Item(4, 9).Value = '基数 xxxx，比例 xx%'
$ws.Cells.Item(4, 10).Value = '13薪还是根据公司业绩提供，是否折扣，折扣比例。'
$ws.Cells.Item(4, 11).Value = '是否打折，比如 xx%。'
$ws.Cells.Item(4, 12).Value = '工位大小，环境，是否提供设备，设备型号种类。'
$ws.Cells.Item(4, 13).Value = '是否有入职就有，是否有前置条件才有。'
$ws.Cells.Item(4, 14).Value = '是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。'
$ws.Cells.Item(4, 15).Value = ''
$ws.Cells.Item(4, 16).Value = ''
$ws.Cells.Item(4, 17).Value = '2022-01-25 07:36:35'

# ---- Row 5 ----
$ws.Cells.Item(5, 2).Value = '南京耀多信息技术有限公司'
$ws.Cells.Item(5, 3).Value = '江苏南京'
$ws.Cells.Item(5, 4).Value = '技术部'
$ws.Cells.Item(5, 5).Value = 'Android'
$ws.Cells.Item(5, 6).Value = '9:00-18:00'
$ws.Cells.Item(5, 7).Value = '1h'
$ws.Cells.Item(5, 8).Value = '一开始996，后来发不起加班费不给加班了，欠的加班费也不发'
$ws.Cells.Item(5, 9).Value = '最低额度'
$ws.Cells.Item(5, 10).Value = '无'
$ws.Cells.Item(5, 11).Value = '八折'
$ws.Cells.Item(5, 12).Value = '提供笔记本'
$ws.Cells.Item(5, 13).Value = '有'
$ws.Cells.Item(5, 14).Value = '钉钉位置打卡'
$ws.Cells.Item(5, 15).Value = '老板阴晴不定，随意开除员工'
$ws.Cells.Item(5, 16).Value = ''
$ws.Cells.Item(5, 17).Value = '2022-01-25 02:22:42'

# ---- Row 6 ----
$ws.Cells.Item(6, 2).Value = '南京希音电子商务有限公司'
$ws.Cells.Item(6, 3).Value = '天溯产业园'
$ws.Cells.Item(6, 4).Value = ''
$ws.Cells.Item(6, 5).Value = '前端'
$ws.Cells.Item(6, 6).Value = '10:00-18:00(到20:00有50补贴)'
$ws.Cells.Item(6, 7).Value = '12:00-13:30'
$ws.Cells.Item(6, 8).Value = '看部门，不强制，有工时排名。'
$ws.Cells.Item(6, 9).Value = '基础工资的8%'
$ws.Cells.Item(6, 10).Value = '看部门盈利情况和个人绩效定'
$ws.Cells.Item(6, 11).Value = '试用期6个月，100%工资不打折'
$ws.Cells.Item(6, 12).Value = '配mac m1+显示器，网吧工作环境，工位挤。'
$ws.Cells.Item(6, 13).Value = '法定年假，可用加班时长来调休'
$ws.Cells.Item(6, 14).Value = '1月3次补卡'
$ws.Cells.Item(6, 15).Value = '抠，舍得给校招生，不舍得给社招生。多余的调休时长换钱200/d'
$ws.Cells.Item(6, 16).Value = ''
$ws.Cells.Item(6, 17).Value = '2022-01-25 01:58:09'

# ---- Row 7 ----
$ws.Cells.Item(7, 2).Value = '慧资环球'
$ws.Cells.Item(7, 3).Value = '白下（年中搬到河西）'
$ws.Cells.Item(7, 4).Value = '研发中心'
$ws.Cells.Item(7, 5).Value = '.NET/Python etc.'
$ws.Cells.Item(7, 6).Value = '自己安排，满8小时工时就好'
$ws.Cells.Item(7, 7).Value = '自己安排'
$ws.Cells.Item(7, 8).Value = '不加班'
$ws.Cells.Item(7, 9).Value = '全额8%'
$ws.Cells.Item(7, 10).Value = '13薪，每年调薪一次'
$ws.Cells.Item(7, 11).Value = '不打折'
$ws.Cells.Item(7, 12).Value = '一个高配台式机或者一个高配Dell工作站笔记本，两个40寸4K显示器 Processor Intel(R) Core(TM) i9-10980XE CPU @ 3.00GHz 3.00 GHz  128GB RAM (新的台式机配置标准)'
$ws.Cells.Item(7, 13).Value = '10 ~ 20天'
$ws.Cells.Item(7, 14).Value = '完全不打卡'
$ws.Cells.Item(7, 15).Value = '内推VX：Just1n'
$ws.Cells.Item(7, 16).Value = ''
$ws.Cells.Item(7, 17).Value = '2022-01-24 14:35:55'

# ---- Row 8 ----
$ws.Cells.Item(8, 2).Value = '零字节'
$ws.Cells.Item(8, 3).Value = '建邺'
$ws.Cells.Item(8, 4).Value = ''
$ws.Cells.Item(8, 5).Value = 'Go/Rust/JS/TS/产品/运营'
$ws.Cells.Item(8, 6).Value = '9：30-6：30'
$ws.Cells.Item(8, 7).Value = '1.5h'
$ws.Cells.Item(8, 8).Value = '不加班'
$ws.Cells.Item(8, 9).Value = '8%'
$ws.Cells.Item(8, 10).Value = '13薪，每年调薪一次'
$ws.Cells.Item(8, 11).Value = '应届生八折，有工作经验的不打折'
$ws.Cells.Item(8, 12).Value = 'macbook pro（入职满三年电脑转赠给员工），每人配一个显示器（24-32寸）'
$ws.Cells.Item(8, 13).Value = '入职转正就享受年假'
$ws.Cells.Item(8, 14).Value = '飞书打卡'
$ws.Cells.Item(8, 15).Value = '节日红包、年度旅游（21年三亚一周）'
$ws.Cells.Item(8, 16).Value = ''
$ws.Cells.Item(8, 17).Value = '2022-01-24 14:32:45'

# ---- Row 9 ----
$ws.Cells.Item(9, 2).Value = '南京力方科技有限公司(力方智充)'
$ws.Cells.Item(9, 3).Value = '雨花台区软件谷科创城'
$ws.Cells.Item(9, 4).Value = '技术部'
$ws.Cells.Item(9, 5).Value = 'Java'
$ws.Cells.Item(9, 6).Value = '9:00-18:00'
$ws.Cells.Item(9, 7).Value = '1.5h'
$ws.Cells.Item(9, 8).Value = '124固定加班到9点'
$ws.Cells.Item(9, 9).Value = '最低，双边合计512'
$ws.Cells.Item(9, 10).Value = '无'
$ws.Cells.Item(9, 11).Value = '三个月，打八折'
$ws.Cells.Item(9, 12).Value = '网吧工位，自带电脑'
$ws.Cells.Item(9, 13).Value = '法定年假'
$ws.Cells.Item(9, 14).Value = '严格打卡，迟打卡扣30，不打卡半天工资'
$ws.Cells.Item(9, 15).Value = ''
$ws.Cells.Item(9, 16).Value = ''
$ws.Cells.Item(9, 17).Value = '2022-01-24 14:29:37'

# ---- Row 10 ----
$ws.Cells.Item(10, 2).Value = '硅基智能'
$ws.Cells.Item(10, 3).Value = '软件大道'
$ws.Cells.Item(10, 4).Value = '创新产品事业群'
$ws.Cells.Item(10, 5).Value = 'Java'
$ws.Cells.Item(10, 6).Value = '9:00-18:30'
$ws.Cells.Item(10, 7).Value = '1.5h'
$ws.Cells.Item(10, 8).Value = '没事到点走，部门氛围卷'
$ws.Cells.Item(10, 9).Value = '基数5500，比例10%'
$ws.Cells.Item(10, 10).Value = '13薪还是根据公司业绩提供，是否折扣，折扣比例。'
$ws.Cells.Item(10, 11).Value = '不打折'
$ws.Cells.Item(10, 12).Value = '网吧工位'
$ws.Cells.Item(10, 13).Value = '满一年才有正常年假，年假次年一月发放（不满一年打折）'
$ws.Cells.Item(10, 14).Value = '是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。'
$ws.Cells.Item(10, 15).Value = ''
$ws.Cells.Item(10, 16).Value = ''
$ws.Cells.Item(10, 17).Value = '2022-01-24 14:25:34'

# ---- Row 11 ----
$ws.Cells.Item(11, 2).Value = '百家云'
$ws.Cells.Item(11, 3).Value = '雨花台软件谷科创城'
$ws.Cells.Item(11, 4).Value = ''
$ws.Cells.Item(11, 5).Value = 'Java'
$ws.Cells.Item(11, 6).Value = '9:00-18:30'
$ws.Cells.Item(11, 7).Value = '1.5h'
$ws.Cells.Item(11, 8).Value = '周1,2,4正常加班，不想加班也行'
$ws.Cells.Item(11, 9).Value = ''
$ws.Cells.Item(11, 10).Value = ''
$ws.Cells.Item(11, 11).Value = '6个月不打折。'
$ws.Cells.Item(11, 12).Value = 'mac笔记本+小米曲面屏显示器'
$ws.Cells.Item(11, 13).Value = '年假次年一月发放，每满一年+1天'
$ws.Cells.Item(11, 14).Value = '每个月有4次迟到补卡机会，早上9.15之前打卡不算迟到'
$ws.Cells.Item(11, 15).Value = ''
$ws.Cells.Item(11, 16).Value = ''
$ws.Cells.Item(11, 17).Value = '2022-01-24 14:21:22'

# ---- Row 12 ----
$ws.Cells.Item(12, 2).Value = '创维南京分公司'
$ws.Cells.Item(12, 3).Value = '雨花云密城'
$ws.Cells.Item(12, 4).Value = 'web后台'
$ws.Cells.Item(12, 5).Value = 'Java'
$ws.Cells.Item(12, 6).Value = '09:30'
$ws.Cells.Item(12, 7).Value = '1.5h'
$ws.Cells.Item(12, 8).Value = '995'
$ws.Cells.Item(12, 9).Value = '工资八折的10%'
$ws.Cells.Item(12, 10).Value = '1个月工资'
$ws.Cells.Item(12, 11).Value = '不打折'
$ws.Cells.Item(12, 12).Value = 'Windows电脑+dell显示器'
$ws.Cells.Item(12, 13).Value = '法定年假'
$ws.Cells.Item(12, 14).Value = '弹性打卡'
$ws.Cells.Item(12, 15).Value = ''
$ws.Cells.Item(12, 16).Value = ''
$ws.Cells.Item(12, 17).Value = '2022-01-24 14:19:34'

# ---- Row 13 ----
$ws.Cells.Item(13, 2).Value = '新视云'
$ws.Cells.Item(13, 3).Value = '雨花台'
$ws.Cells.Item(13, 4).Value = ''
$ws.Cells.Item(13, 5).Value = 'Java'
$ws.Cells.Item(13, 6).Value = '9:00-17:30'
$ws.Cells.Item(13, 7).Value = '1h'
$ws.Cells.Item(13, 8).Value = '看部门，业务部门偶尔加班，技术支持部门基本不加班'
$ws.Cells.Item(13, 9).Value = '基数5k，比例8%'
$ws.Cells.Item(13, 10).Value = '固定13薪'
$ws.Cells.Item(13, 11).Value = '3年合同，试用期总共6个月，前三个月8折，后三个月全薪'
$ws.Cells.Item(13, 12).Value = '配笔记本+显示器'
$ws.Cells.Item(13, 13).Value = '5天年假+5天带薪病假（入职自动折算当年年假）'
$ws.Cells.Item(13, 14).Value = '不打卡'
$ws.Cells.Item(13, 15).Value = ''
$ws.Cells.Item(13, 16).Value = ''
$ws.Cells.Item(13, 17).Value = '2022-01-24 14:17:01'

# ---- Row 14 ----
$ws.Cells.Item(14, 2).Value = '华为'
$ws.Cells.Item(14, 3).Value = '华为南研所'
$ws.Cells.Item(14, 4).Value = ''
$ws.Cells.Item(14, 5).Value = 'Java'
$ws.Cells.Item(14, 6).Value = '9:00'
$ws.Cells.Item(14, 7).Value = '12:00-13:40'
$ws.Cells.Item(14, 8).Value = '看部门情况。好部门：124加班8：30，35正常下班,差部门：天天11点以后'
$ws.Cells.Item(14, 9).Value = '基础工资的5%'
$ws.Cells.Item(14, 10).Value = '看部门盈利情况和个人绩效定'
$ws.Cells.Item(14, 11).Value = '试用期6个月，100%工资不打折'
$ws.Cells.Item(14, 12).Value = '配win台式机+双屏'
$ws.Cells.Item(14, 13).Value = '没签奋斗协议的5天，但一般不给休，第二年可以换成钱。签了的自愿放弃年假了'
$ws.Cells.Item(14, 14).Value = '必须按时打卡'
$ws.Cells.Item(14, 15).Value = ''
$ws.Cells.Item(14, 16).Value = ''
$ws.Cells.Item(14, 17).Value = '2022-01-24 14:17:32'

# ---- Row 15 ----
$ws.Cells.Item(15, 2).Value = '满帮'
$ws.Cells.Item(15, 3).Value = '雨花区万博科技园'
$ws.Cells.Item(15, 4).Value = ''
$ws.Cells.Item(15, 5).Value = 'Java'
$ws.Cells.Item(15, 6).Value = '9:00-18:30'
$ws.Cells.Item(15, 7).Value = '1.5h'
$ws.Cells.Item(15, 8).Value = '看部门，不强制， 周五基本不加，还有每月一天奋斗日（年底算工资）， 据说要取消了'
$ws.Cells.Item(15, 9).Value = '全额8%'
$ws.Cells.Item(15, 10).Value = '上下半年绩效'
$ws.Cells.Item(15, 11).Value = ''
$ws.Cells.Item(15, 12).Value = '联想'
$ws.Cells.Item(15, 13).Value = ''
$ws.Cells.Item(15, 14).Value = ''
$ws.Cells.Item(15, 15).Value = ''
$ws.Cells.Item(15, 16).Value = ''
$ws.Cells.Item(15, 17).Value = '2022-01-24 14:10:47'

# ---- Row 16 ----
$ws.Cells.Item(16, 2).Value = 'A示例xxx公司'
$ws.Cells.Item(16, 3).Value = 'xx区'
$ws.Cells.Item(16, 4).Value = 'xxx事业部'
$ws.Cells.Item(16, 5).Value = 'Java'
$ws.Cells.Item(16, 6).Value = '9:00-18:30'
$ws.Cells.Item(16, 7).Value = '1.5h'
$ws.Cells.Item(16, 8).Value = '135 加班，24 正常；大小周等等'
$ws.Cells.Item(16, 9).Value = '基数 xxxx，比例 xx%'
$ws.Cells.Item(16, 10).Value = '13薪还是根据公司业绩提供，是否折扣，折扣比例。'
$ws.Cells.Item(16, 11).Value = '是否打折，比如 xx%。'
$ws.Cells.Item(16, 12).Value = '工位大小，环境，是否提供设备，设备型号种类。'
$ws.Cells.Item(16, 13).Value = '是否有入职就有，是否有前置条件才有。'
$ws.Cells.Item(16, 14).Value = '是否严格打卡，使用的软件或者方式（比如钉钉或人脸识别）。'
$ws.Cells.Item(16, 15).Value = ''
$ws.Cells.Item(16, 16).Value = ''
$ws.Cells.Item(16, 17).Value = '2022-01-24 13:11:01'

# ---- New row 16 index cell (column A), continuing the sequence 0..14 ----
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)

$excel.CutCopyMode = $false
$ws.Range("A1").Select()